$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.813.67'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.92%  '
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.506.71'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.72%  '
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.69%  '
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.77'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.81%  '
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.499.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.57%  '
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.611'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.37%  '
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.622'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.38%  '
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.79%  '
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.65%  '
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.16'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.91%  '
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.071.33'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.65%  '
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '649.86'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.51%  '
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.881.64'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.508.61'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.23%  '
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.88%  '
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.38'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.85%  '
$ws.Range("E20").ClearFormats()

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.952'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.10%  '
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.09'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.71%  '
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.35'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.44'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.37%  '
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -7.06%  '
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.95%  '
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.45%  '
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.39'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.66%  '
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.73'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.69%  '
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.24'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -7.38%  '
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.75'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.25%  '
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.65'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.19%  '
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '579.39'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +10.71%  '
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '61.41'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.35%  '
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.753.33'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.54%  '
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0795'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -9.39%  '
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.61'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("B41").Value = 'CoreDAO'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.79'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +45.55%  '
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.38%  '
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.375'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.44%  '
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.08%  '
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '34.25'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.92%  '
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.44%  '
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.70%  '
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.82%  '
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.44%  '
$ws.Range("E51").ClearFormats()
